# Append the 2025/11/14 mod-count data point as a new row (row 5) to the
# "ModCounts" sheet, following the same pattern as the existing rows 3/4
# (date, game name, numeric mod count; all centered like rows 3:4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new values first. The date column is entered with a leading
# apostrophe so it is stored as literal text ("2025/11/14"), matching the
# existing Date column cells, instead of being auto-parsed into a date
# serial number by Excel's smart input handling.
$ws.Range("A5").Value = "'2025/11/14"
$ws.Range("B5").Value = "逃离鸭科夫"
$ws.Range("C5").Value = 1089

# Copy the formatting (centered alignment) from the previous data row so
# the new row matches rows 3:4 visually/stylistically.
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
